# Update cryptos list figures (prices & 1h volume change)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.614.73"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.27%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.865.04"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.50%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "326.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.27%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.004"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.23%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4620"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.54%  "
$ws.Range("E8").Value = "  +1.24%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07907"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.63%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9683"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.25%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "22.28"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.81%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.886.06"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.99%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.715"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.07%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.928"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.23%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06949"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.20%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.19"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.66%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.005"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.26%  "
$ws.Range("E18").Value = "  +0.86%  "
$ws.Range("E19").Value = "  +0.28%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.004"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.33%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "28.638.32"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.22%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.307"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.53%  "
$ws.Range("E23").Value = "  +0.84%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.123"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.19%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.174.68"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.13%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "153.72"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.31"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.48%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.708"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.991"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.53%  "
$ws.Range("E30").Value = "  +2.01%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09348"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.47%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9303"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.16%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.323"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.78%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.339"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.14%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.356"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.59%  "
$ws.Range("E36").Value = "  -3.19%  "
$ws.Range("E37").Value = "  -1.56%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.151"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.897"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.94%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5640"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.37%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "9.906"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.99%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1778"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.19%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.07230"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.90%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "11.68"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.63%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5313"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.38%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.150"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.36%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.136"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -8.02%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.843"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.34%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "113.23"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.38%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.004"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.34%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.344"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.68%  "
